$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2025-10-28 03:53:53"
$ws.Range("B5").Value = "'2025-10-27"
$ws.Range("C5").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("D5").Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/ListPrice27102025.pdf"
